# Applies two kinds of edits described by the diff:
#   1. Remove the "Heading2" paragraph style from the five section
#      headings (Introduction, Virtue Framework, Utilitarian Framework,
#      Deontological Framework, Conclusion) so the <w:pPr>/<w:pStyle>
#      element disappears entirely (not just switched to "Normal").
#   2. Replace the inline author-citation placeholders with the new
#      "Ref-xxxxxxx" reference keys used by the citation checker.

$d = $word.ActiveDocument

function Remove-HeadingParagraphStyle($doc, $headingText) {
    $target = $null
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $candidate = $doc.Paragraphs.Item($i)
        if ($candidate.Range.Text.TrimEnd([char]13) -eq $headingText) {
            $target = $candidate
            break
        }
    }
    if ($target -eq $null) {
        Write-Host "heading not found: $headingText"
        return
    }
    # Deleting the paragraph's whole range (text + paragraph mark) drops
    # the w:pPr/w:pStyle that lives on that mark; re-inserting the text
    # with a fresh paragraph mark creates a <w:p> with no pPr at all,
    # matching the diff (style removed, not just reset to "Normal").
    $insertPos = $target.Range.Start
    $fullRange = $doc.Range($target.Range.Start, $target.Range.End)
    $fullRange.Delete()
    $newRange = $doc.Range($insertPos, $insertPos)
    $newRange.InsertBefore("$headingText`r")
}

Remove-HeadingParagraphStyle $d "Introduction"
Remove-HeadingParagraphStyle $d "Virtue Framework"
Remove-HeadingParagraphStyle $d "Utilitarian Framework"
Remove-HeadingParagraphStyle $d "Deontological Framework"
Remove-HeadingParagraphStyle $d "Conclusion"

# Citation placeholder -> Ref-ID replacements. Plain author mentions
# outside parentheses (e.g. "According to Königs, ..." or "Sætra
# highlights ...") are left untouched; only the parenthetical citation
# markers change.
$d.Content.Find.Execute("(Russell)", $true, $false, $false, $false, $false, $true, 1, $false, "(Ref-f399966)", 2) | Out-Null
$d.Content.Find.Execute("(Bellazzi and Boyneburgk)", $true, $false, $false, $false, $false, $true, 1, $false, "(Ref-f399966)", 2) | Out-Null

$d.Content.Find.Execute("(Alexander)", $true, $false, $false, $false, $false, $true, 1, $false, "(Ref-f156724)", 2) | Out-Null
$d.Content.Find.Execute("value of privacy and liberty (Sætra)", $true, $false, $false, $false, $false, $true, 1, $false, "value of privacy and liberty (Ref-f156724)", 2) | Out-Null

$d.Content.Find.Execute("(Königs)", $true, $false, $false, $false, $false, $true, 1, $false, "(Ref-f651696)", 2) | Out-Null
$d.Content.Find.Execute("both privacy and liberty (Sætra)", $true, $false, $false, $false, $false, $true, 1, $false, "both privacy and liberty (Ref-f651696)", 2) | Out-Null

Write-Host "done"
